# Update the date line at the top of the document.
$d = $word.ActiveDocument

$d.Content.Find.Execute("2024-06-27 Thursday", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "2024-06-28 Friday", 2)

# Update the practice-problem table. Each data row in the table is
# followed by three blank spacer rows, so the populated rows are
# Word row indices 1, 5, 9, 13, 17 (1-based), each with 5 columns.
# We address cells by position (row, column) rather than doing a
# document-wide text replace, because several of the old values repeat
# (e.g. "63÷8=7, 7" appears twice) but map to different new values.

$tbl = $d.Tables.Item(1)

$tbl.Cell(1, 1).Range.Text  = "33÷9=3, 6"
$tbl.Cell(1, 2).Range.Text  = "58÷7=8, 2"
$tbl.Cell(1, 3).Range.Text  = "31÷2=15, 1"
$tbl.Cell(1, 4).Range.Text  = "56÷4=14, 0"
$tbl.Cell(1, 5).Range.Text  = "51÷3=17, 0"

$tbl.Cell(5, 1).Range.Text  = "58÷4=14, 2"
$tbl.Cell(5, 2).Range.Text  = "33÷9=3, 6"
$tbl.Cell(5, 3).Range.Text  = "31÷9=3, 4"
$tbl.Cell(5, 4).Range.Text  = "64÷5=12, 4"
$tbl.Cell(5, 5).Range.Text  = "32÷6=5, 2"

$tbl.Cell(9, 1).Range.Text  = "81÷6=13, 3"
$tbl.Cell(9, 2).Range.Text  = "28÷3=9, 1"
$tbl.Cell(9, 3).Range.Text  = "84÷4=21, 0"
$tbl.Cell(9, 4).Range.Text  = "91÷4=22, 3"
$tbl.Cell(9, 5).Range.Text  = "27÷5=5, 2"

$tbl.Cell(13, 1).Range.Text = "67÷4=16, 3"
$tbl.Cell(13, 2).Range.Text = "31÷4=7, 3"
$tbl.Cell(13, 3).Range.Text = "44÷2=22, 0"
$tbl.Cell(13, 4).Range.Text = "66÷5=13, 1"
$tbl.Cell(13, 5).Range.Text = "56÷5=11, 1"

$tbl.Cell(17, 1).Range.Text = "97÷7=13, 6"
$tbl.Cell(17, 2).Range.Text = "56÷2=28, 0"
$tbl.Cell(17, 3).Range.Text = "86÷6=14, 2"
$tbl.Cell(17, 4).Range.Text = "57÷5=11, 2"
$tbl.Cell(17, 5).Range.Text = "56÷4=14, 0"

Write-Output "edit applied"
